# Weekly update: insert a new weekly record (row) for "Camote" at Vega Central
# Mapocho de Santiago. This pushes the existing rows 43-51 down to 44-52, and
# fills the newly inserted row 43 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 43 (existing rows 43:51 shift to 44:52)
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new weekly record.
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 44508
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100114002
$ws.Cells.Item(43, 7).Value = "Camote"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 1420
$ws.Cells.Item(43, 11).Value = 13000
$ws.Cells.Item(43, 12).Value = 14000
$ws.Cells.Item(43, 13).Value = 13500
$ws.Cells.Item(43, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(43, 15).Value = "Perú"
$ws.Cells.Item(43, 16).Value = 750
$ws.Cells.Item(43, 17).Value = 18
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date style/format as the other
# rows in column D (style index used by D2:D52).
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
